$d = $word.ActiveDocument

$replacements = @(
    @("638×5=", "546×4="),
    @("114×6=", "807×3="),
    @("687×9=", "607×8="),
    @("534×9=", "783×2="),
    @("629×2=", "429×5="),
    @("797×3=", "591×4="),
    @("119×3=", "361×3="),
    @("146×7=", "458×2="),
    @("493×2=", "351×9="),
    @("301×2=", "220×9="),
    @("989×3=", "567×5="),
    @("743×4=", "662×3="),
    @("397×2=", "293×2="),
    @("926×4=", "400×8="),
    @("891×6=", "954×9="),
    @("919×5=", "978×3="),
    @("515×6=", "309×3="),
    @("932×6=", "921×8="),
    @("180×3=", "519×5="),
    @("612×2=", "153×7="),
    @("996×4=", "179×4="),
    @("806×2=", "564×6="),
    @("560×9=", "825×9="),
    @("102×2=", "489×6="),
    @("518×5=", "611×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
